# Update "DEG Number of Identified Genes.xlsx":
#  - Insert a new first worksheet "M_MUT_and_WT_M_E18_WB" (an E18 whole-brain
#    DEG sheet, still unfilled / all data cleared) ahead of the three
#    existing P30/P60/P120 CORT sheets.
#  - Populate it with the same header row / cluster-name column as the other
#    sheets, wire up the same roll-up formulas in B:D, and leave the raw
#    per-method counts (E:K) blank, matching the other sheets' layout.
#  - Leave the other three sheets' data untouched (only their selection /
#    active-sheet bookkeeping changes as a side effect of which sheet is
#    active when the workbook is saved).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Create the new sheet and move it to the front of the tab strip.
# ---------------------------------------------------------------------------
$new = $wb.Worksheets.Add()
$new.Name = "M_MUT_and_WT_M_E18_WB"
$new.Move($wb.Worksheets.Item(1))

$ws = $wb.Worksheets.Item(1)
$p30 = $wb.Worksheets.Item(2)   # M_MUT_and_WT_M_P30_CORT  (template for layout)

# ---------------------------------------------------------------------------
# 2. Column widths -- match the other sheets (B:K only, A keeps default).
# ---------------------------------------------------------------------------
for ($i = 2; $i -le 11; $i++) {
    $ws.Columns.Item($i).ColumnWidth = $p30.Columns.Item($i).ColumnWidth
}

# ---------------------------------------------------------------------------
# 3. Header row (row 1) -- identical text/order to the other three sheets.
# ---------------------------------------------------------------------------
$headers = @("Cluster", "DESeq2 (tot)", "Limma (tot)", "EdgeR (tot)", `
             "DESeq2 Only", "Limma Only", "EdgeR Only", `
             "DESeq2 & Limma", "DESeq2 & EdgeR", "Limma & EdgeR", "All Methods")
$ws.Range("A1:K1").Value = $headers
$ws.Range("A1:K1").Font.Bold = $true
$ws.Range("A1:K1").HorizontalAlignment = -4108
$ws.Range("A1:K1").VerticalAlignment = -4108

# ---------------------------------------------------------------------------
# 4. Cluster-name column (A2:A15) -- same 14 clusters, same order.
# ---------------------------------------------------------------------------
$clusters = @("L2_3_IT", "L6", "Sst", "L5", "L4", "Pvalb", "Sncg", `
              "Non-neuronal", "Oligo", "Vip", "Lamp5", "Astro", "Peri", "Endo")
for ($r = 0; $r -lt $clusters.Length; $r++) {
    $row = $r + 2
    $ws.Cells.Item($row, 1).Value = $clusters[$r]
}
$ws.Range("A2:A15").HorizontalAlignment = -4108
$ws.Range("A2:A15").VerticalAlignment = -4108

# ---------------------------------------------------------------------------
# 5. Roll-up formulas in B:D (same shape as the other sheets); E:K stay
#    empty -- this tab has not been filled in with per-method DEG counts
#    yet, so every total comes out to 0.
# ---------------------------------------------------------------------------
for ($r = 2; $r -le 15; $r++) {
    $ws.Cells.Item($r, 2).Formula = "=E$r+H$r+I$r+K$r"
    $ws.Cells.Item($r, 3).Formula = "=F$r+H$r+J$r+K$r"
    $ws.Cells.Item($r, 4).Formula = "=G$r+J$r+I$r+K$r"
}
$ws.Range("B2:K15").HorizontalAlignment = -4108
$ws.Range("B2:K15").VerticalAlignment = -4108

# ---------------------------------------------------------------------------
# 6. Trailing blank rows (16:18, columns B:K) -- present on the sheet but
#    with no content/centring, matching the spacer rows under the table.
# ---------------------------------------------------------------------------
$ws.Range("B16:K18").HorizontalAlignment = -4142
$ws.Range("B16:K18").Value = ""

$ws.Range("A1").Select()

Write-Output "Added sheet M_MUT_and_WT_M_E18_WB as the first tab."
